# Apply the "agregar codigo de barras" patient-data update.
# Updates the patient-record fields on the active sheet and clears the
# now-unused "Tipo de Consulta" free-text cell (D16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1er Apellido / 2do Apellido / 1er Nombre / 2do Nombre + No. Expediente
$ws.Range("A6").Value = "TORRES  CHAN  JESSICA  MARIBEL"
$ws.Range("G6").Value = "/201773345"

# Fecha de Nacimiento / Edad / Lugar de Nacimiento
$ws.Range("A9").Value = "1988-12-09"
$ws.Range("D9").Value = "28 AÑOS"
$ws.Range("E9").Value = "GUATEMALA"

# Estado Civil / Ocupacion / Documento de Identificacion
$ws.Range("A11").Value = "SOLTERO"
$ws.Range("C11").Value = "MEDICO"
$ws.Range("G11").Value = ""

# Emergencia: Nombre / Parentesco / Direccion / Telefono
$ws.Range("A13").Value = "ALEXANDER TORRES"
$ws.Range("D13").Value = "HERMANO"
$ws.Range("E13").Value = "4A CALLE A 15-27 MIXCO"
$ws.Range("G13").Value = "59591912"

# Fecha de la asistencia Medica: Hora / Area de urgencia / Fecha
$ws.Range("D14").Value = "Hora: 10:30:14"
$ws.Range("E14").Value = "Area de urgencia: null"
$ws.Range("A15").Value = "20/11/2017"

# Tipo de Consulta free-text cell no longer used
$ws.Range("D16").Value = ""
